$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above current row 9 ("un_franzosa_ControlvsCD_Fp")
# so it becomes "un_franzosa_ControlvsCD_ConvCD" with the same values as
# the neighboring "un_franzosa_ControlvsCD_Age"/"un_franzosa_ControlvsCD_Fp" rows.
$ws.Rows.Item(9).Insert()
$ws.Range("A9").Value = "un_franzosa_ControlvsCD_ConvCD"
$ws.Range("B9").Value = 0
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1
$ws.Range("H9").Value = 1

# Insert a new row above row 15 ("un_franzosa_ControlvsUC_Fp", after the
# previous insertion shifted it down from row 14) so it becomes
# "un_franzosa_ControlvsUC_ConvUC".
$ws.Rows.Item(15).Insert()
$ws.Range("A15").Value = "un_franzosa_ControlvsUC_ConvUC"
$ws.Range("B15").Value = 0
$ws.Range("C15").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 0
